# Applies the "読み取り元" authoring pass:
#  - renames the sheet from "Sheet1" to "読み取り元"
#  - clears the leftover placeholder formatting that used to fill B:P on
#    rows 2-4 (and C5), removing the now-empty cells entirely
#  - strips the (hidden/white-text) comment styling from the cells that
#    still hold real text, so they fall back to the default cell style
#  - re-applies just the white font (no fill) to G2, which is the one
#    "placeholder" cell that keeps a non-default style
#  - drops the stray selection left over on C7

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- rename the sheet ------------------------------------------------
$ws.Name = "読み取り元"

# --- cells that were only carrying leftover blank formatting: clear
#     formatting AND contents so the cell disappears entirely ----------
# (kept as single contiguous blocks - comma/multi-area Ranges only
# clear their first area in this runtime, so each block is separate)
$emptyRanges = @(
    "B2:B2",
    "D2:F2",
    "H2:P2",
    "B3:B3",
    "D3:P3",
    "B4:B4",
    "D4:P4",
    "C5"
)
foreach ($addr in $emptyRanges) {
    $r = $ws.Range($addr)
    $r.ClearFormats()
    $r.ClearContents()
}

# --- cells that still hold real text: drop the special formatting but
#     keep the text itself -------------------------------------------
$textCells = @("C2", "C3", "C4", "C7", "C8", "C44", "C58")
foreach ($addr in $textCells) {
    $ws.Range($addr).ClearFormats()
}

# --- G2 keeps a non-default style: white font text, no fill -----------
$g2 = $ws.Range("G2")
$g2.ClearFormats()
$g2.Font.ThemeColor = 2

# --- clear the leftover selection rectangle on C7 ----------------------
$ws.Range("A1").Select()
